$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1499
$ws.Range("F3").Value = 829
$ws.Range("F5").Value = 881
$ws.Range("F6").Value = 492
$ws.Range("F7").Value = 7347
$ws.Range("F11").Value = 5361
$ws.Range("F12").Value = 555
$ws.Range("F14").Value = 7411
$ws.Range("F15").Value = 8705
$ws.Range("F16").Value = 162
$ws.Range("F17").Value = 1131
$ws.Range("F18").Value = 858
$ws.Range("F19").Value = 4360
$ws.Range("F20").Value = 657
$ws.Range("F21").Value = 188
$ws.Range("F22").Value = 81
$ws.Range("F25").Value = 1177
$ws.Range("F26").Value = 85
$ws.Range("F27").Value = 1617
$ws.Range("F29").Value = 875
$ws.Range("F30").Value = 1845
$ws.Range("F31").Value = 316
$ws.Range("F32").Value = 2234
$ws.Range("F34").Value = 98
$ws.Range("F35").Value = 1415
$ws.Range("F40").Value = 2908
$ws.Range("F41").Value = 4005
$ws.Range("F42").Value = 186
$ws.Range("F43").Value = 36
$ws.Range("F44").Value = 404
$ws.Range("F45").Value = 504
$ws.Range("F48").Value = 157

$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 3
$ws.Range("F19").Value = 2

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5087

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 5087
$ws.Range("F3").Value = 1499
$ws.Range("F4").Value = 829
$ws.Range("F6").Value = 881
$ws.Range("F7").Value = 492
$ws.Range("F8").Value = 3
$ws.Range("F10").Value = 5361
$ws.Range("F11").Value = 555
$ws.Range("F12").Value = 7411
$ws.Range("F14").Value = 162
$ws.Range("F15").Value = 1131
$ws.Range("F16").Value = 858
$ws.Range("F17").Value = 4360
$ws.Range("F18").Value = 657
$ws.Range("F19").Value = 188
$ws.Range("F23").Value = 1177
$ws.Range("F24").Value = 85
$ws.Range("F25").Value = 1617
$ws.Range("F27").Value = 875
$ws.Range("F28").Value = 1845
$ws.Range("F29").Value = 316
$ws.Range("F30").Value = 2234
$ws.Range("F40").Value = 4005
$ws.Range("F42").Value = 186
$ws.Range("F43").Value = 36
$ws.Range("F44").Value = 404
$ws.Range("F45").Value = 504
$ws.Range("F47").Value = 157
